$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("H2").Value = 98
$ws.Range("I2").Value = 229
$ws.Range("J2").Value = 894
$ws.Range("K2").Value = 4
$ws.Range("L2").Value = 273
$ws.Range("M2").Value = 14
$ws.Range("N2").Value = 187
$ws.Range("P2").Value = 5
$ws.Range("Q2").Value = 2
$ws.Range("R2").Value = 11
$ws.Range("S2").Value = 99
$ws.Range("T2").Value = 167
$ws.Range("U2").Value = 9
$ws.Range("V2").Value = 1425
$ws.Range("X2").Value = 1471
$ws.Range("Y2").Value = 2
$ws.Range("Z2").Value = 21
$ws.Range("AA2").Value = 12
